# Update TPM-derived NATMI metrics for Sema3a-Plxna4 ligand-receptor pairs
# (rows 2-16) with the refreshed scripted output.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.458056666666667
$ws.Range("H2").Value = 4.374169999999999
$ws.Range("I2").Value = 0.2323568509805328
$ws.Range("J2").Value = 0.2323568509805327
$ws.Range("M2").Value = 6.602366333333333
$ws.Range("N2").Value = 19.807099
$ws.Range("O2").Value = 0.7068089336605662
$ws.Range("P2").Value = 0.7068089336605662
$ws.Range("Q2").Value = 9.626624248092222
$ws.Range("R2").Value = 86.63961823282999
$ws.Range("S2").Value = 0.1642318980702774
$ws.Range("T2").Value = 0.1642318980702774
# Row 3
$ws.Range("G3").Value = 1.458056666666667
$ws.Range("H3").Value = 4.374169999999999
$ws.Range("I3").Value = 0.2323568509805328
$ws.Range("J3").Value = 0.2323568509805327
$ws.Range("O3").Value = 0.2539225612198319
$ws.Range("P3").Value = 0.2539225612198319
$ws.Range("Q3").Value = 3.458384534441111
$ws.Range("R3").Value = 31.12546080997
$ws.Range("S3").Value = 0.05900064671795168
$ws.Range("T3").Value = 0.05900064671795167
# Row 4
$ws.Range("G4").Value = 1.458056666666667
$ws.Range("H4").Value = 4.374169999999999
$ws.Range("I4").Value = 0.2323568509805328
$ws.Range("J4").Value = 0.2323568509805327
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.04677800000000001
$ws.Range("N4").Value = 0.140334
$ws.Range("O4").Value = 0.005007766402153183
$ws.Range("P4").Value = 0.005007766402153183
$ws.Range("Q4").Value = 0.06820497475333334
$ws.Range("R4").Value = 0.6138447727799999
$ws.Range("S4").Value = 0.001163588831650426
$ws.Range("T4").Value = 0.001163588831650426
# Row 5
$ws.Range("G5").Value = 1.458056666666667
$ws.Range("H5").Value = 4.374169999999999
$ws.Range("I5").Value = 0.2323568509805328
$ws.Range("J5").Value = 0.2323568509805327
$ws.Range("M5").Value = 0.112148
$ws.Range("N5").Value = 0.336444
$ws.Range("O5").Value = 0.01200587854266269
$ws.Range("P5").Value = 0.01200587854266268
$ws.Range("Q5").Value = 0.1635181390533333
$ws.Range("R5").Value = 1.47166325148
$ws.Range("S5").Value = 0.00278964813142785
$ws.Range("T5").Value = 0.002789648131427849
# Row 6
$ws.Range("G6").Value = 1.458056666666667
$ws.Range("H6").Value = 4.374169999999999
$ws.Range("I6").Value = 0.2323568509805328
$ws.Range("J6").Value = 0.2323568509805327
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.2078846666666666
$ws.Range("N6").Value = 0.6236539999999999
$ws.Range("O6").Value = 0.02225486017478616
$ws.Range("P6").Value = 0.02225486017478616
$ws.Range("Q6").Value = 0.303107624131111
$ws.Range("R6").Value = 2.727968617179999
$ws.Range("S6").Value = 0.00517106922922538
$ws.Range("T6").Value = 0.005171069229225379
# Row 7
$ws.Range("H7").Value = 5.708772
$ws.Range("I7").Value = 0.3032511962008422
$ws.Range("J7").Value = 0.3032511962008422
$ws.Range("M7").Value = 6.602366333333333
$ws.Range("N7").Value = 19.807099
$ws.Range("O7").Value = 0.7068089336605662
$ws.Range("P7").Value = 0.7068089336605662
$ws.Range("Q7").Value = 12.563801352492
$ws.Range("R7").Value = 113.074212172428
$ws.Range("S7").Value = 0.2143406546180084
$ws.Range("T7").Value = 0.2143406546180084
# Row 8
$ws.Range("H8").Value = 5.708772
$ws.Range("I8").Value = 0.3032511962008422
$ws.Range("J8").Value = 0.3032511962008422
$ws.Range("O8").Value = 0.2539225612198319
$ws.Range("P8").Value = 0.2539225612198319
$ws.Range("Q8").Value = 4.513571442227999
$ws.Range("S8").Value = 0.0770023204322956
$ws.Range("T8").Value = 0.0770023204322956
# Row 9
$ws.Range("H9").Value = 5.708772
$ws.Range("I9").Value = 0.3032511962008422
$ws.Range("J9").Value = 0.3032511962008422
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.04677800000000001
$ws.Range("N9").Value = 0.140334
$ws.Range("O9").Value = 0.005007766402153183
$ws.Range("P9").Value = 0.005007766402153183
$ws.Range("Q9").Value = 0.08901497887200001
$ws.Range("R9").Value = 0.801134809848
$ws.Range("S9").Value = 0.001518611151747341
$ws.Range("T9").Value = 0.001518611151747341
# Row 10
$ws.Range("H10").Value = 5.708772
$ws.Range("I10").Value = 0.3032511962008422
$ws.Range("J10").Value = 0.3032511962008422
$ws.Range("M10").Value = 0.112148
$ws.Range("N10").Value = 0.336444
$ws.Range("O10").Value = 0.01200587854266269
$ws.Range("P10").Value = 0.01200587854266268
$ws.Range("Q10").Value = 0.213409120752
$ws.Range("R10").Value = 1.920682086768
$ws.Range("S10").Value = 0.003640797029504484
$ws.Range("T10").Value = 0.003640797029504484
# Row 11
$ws.Range("H11").Value = 5.708772
$ws.Range("I11").Value = 0.3032511962008422
$ws.Range("J11").Value = 0.3032511962008422
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.2078846666666666
$ws.Range("N11").Value = 0.6236539999999999
$ws.Range("O11").Value = 0.02225486017478616
$ws.Range("P11").Value = 0.02225486017478616
$ws.Range("Q11").Value = 0.3955887214319999
$ws.Range("R11").Value = 3.560298492887999
$ws.Range("S11").Value = 0.006748812969286386
$ws.Range("T11").Value = 0.006748812969286386
# Row 12
$ws.Range("G12").Value = 2.914094333333333
$ws.Range("H12").Value = 8.742283
$ws.Range("I12").Value = 0.4643919528186251
$ws.Range("J12").Value = 0.4643919528186251
$ws.Range("M12").Value = 6.602366333333333
$ws.Range("N12").Value = 19.807099
$ws.Range("O12").Value = 0.7068089336605662
$ws.Range("P12").Value = 0.7068089336605662
$ws.Range("Q12").Value = 19.23991831855744
$ws.Range("R12").Value = 173.159264867017
$ws.Range("S12").Value = 0.3282363809722804
$ws.Range("T12").Value = 0.3282363809722804
# Row 13
$ws.Range("G13").Value = 2.914094333333333
$ws.Range("H13").Value = 8.742283
$ws.Range("I13").Value = 0.4643919528186251
$ws.Range("J13").Value = 0.4643919528186251
$ws.Range("O13").Value = 0.2539225612198319
$ws.Range("P13").Value = 0.2539225612198319
$ws.Range("Q13").Value = 6.911980175189222
$ws.Range("R13").Value = 62.207821576703
$ws.Range("S13").Value = 0.1179195940695846
$ws.Range("T13").Value = 0.1179195940695846
# Row 14
$ws.Range("G14").Value = 2.914094333333333
$ws.Range("H14").Value = 8.742283
$ws.Range("I14").Value = 0.4643919528186251
$ws.Range("J14").Value = 0.4643919528186251
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 0.04677800000000001
$ws.Range("N14").Value = 0.140334
$ws.Range("O14").Value = 0.005007766402153183
$ws.Range("P14").Value = 0.005007766402153183
$ws.Range("Q14").Value = 0.1363155047246667
$ws.Range("R14").Value = 1.226839542522
$ws.Range("S14").Value = 0.002325566418755416
$ws.Range("T14").Value = 0.002325566418755416
# Row 15
$ws.Range("G15").Value = 2.914094333333333
$ws.Range("H15").Value = 8.742283
$ws.Range("I15").Value = 0.4643919528186251
$ws.Range("J15").Value = 0.4643919528186251
$ws.Range("M15").Value = 0.112148
$ws.Range("N15").Value = 0.336444
$ws.Range("O15").Value = 0.01200587854266269
$ws.Range("P15").Value = 0.01200587854266268
$ws.Range("Q15").Value = 0.3268098512946667
$ws.Range("R15").Value = 2.941288661652
$ws.Range("S15").Value = 0.005575433381730353
$ws.Range("T15").Value = 0.005575433381730352
# Row 16
$ws.Range("G16").Value = 2.914094333333333
$ws.Range("H16").Value = 8.742283
$ws.Range("I16").Value = 0.4643919528186251
$ws.Range("J16").Value = 0.4643919528186251
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.2078846666666666
$ws.Range("N16").Value = 0.6236539999999999
$ws.Range("O16").Value = 0.02225486017478616
$ws.Range("P16").Value = 0.02225486017478616
$ws.Range("Q16").Value = 0.6057955291202222
$ws.Range("R16").Value = 5.452159762081999
$ws.Range("S16").Value = 0.01033497797627439
$ws.Range("T16").Value = 0.01033497797627439

Write-Host "Updated 180 cells with refreshed TPM values"
